# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 978.45
$ws.Range("J58").Value = 1581.75
$ws.Range("L58").Value = 4745.25
$ws.Range("N58").Value = -5045.25

# Row 70
$ws.Range("H70").Value = 3349.5833
$ws.Range("J70").Value = 5085.7144
$ws.Range("L70").Value = 15257.1432
$ws.Range("N70").Value = -15797.1432

# Row 73
$ws.Range("H73").Value = 3349.5833
$ws.Range("J73").Value = 5085.7144
$ws.Range("L73").Value = 15257.1432
$ws.Range("N73").Value = -17129.1432

# Row 88
$ws.Range("H88").Value = 2187
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 874
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 874
$ws.Range("M88").Value = -3094
$ws.Range("N88").Value = -1686

# Row 91
$ws.Range("H91").Value = 2187
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 874
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 874
$ws.Range("M91").Value = -2096
$ws.Range("N91").Value = -3682

# Row 113
$ws.Range("H113").Value = 1336.2727
$ws.Range("I113").Value = 1425
$ws.Range("J113").Value = 449
$ws.Range("K113").Value = 1425
$ws.Range("L113").Value = 449
$ws.Range("M113").Value = 1829
$ws.Range("N113").Value = -6957

# Row 116
$ws.Range("H116").Value = 4585.5713
$ws.Range("J116").Value = 4299.5
$ws.Range("L116").Value = 4299.5
$ws.Range("N116").Value = -11183.5

# Row 121
$ws.Range("H121").Value = 1111.6666
$ws.Range("J121").Value = 1111.6666
$ws.Range("L121").Value = 3334.9998
$ws.Range("N121").Value = -6828.9998

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 138
$ws.Range("H138").Value = 2327.1667
$ws.Range("J138").Value = 2656
$ws.Range("L138").Value = 7968
$ws.Range("N138").Value = -18248

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3624.2273
$ws.Range("I61").Value = 2631.4119
$ws.Range("K61").Value = 2631.4119
$ws.Range("M61").Value = -2419.4119

# Row 74
$ws.Range("H74").Value = 1837.875
$ws.Range("I74").Value = 1800.4286
$ws.Range("K74").Value = 1800.4286
$ws.Range("M74").Value = -926.4286

# Row 77
$ws.Range("H77").Value = 1837.875
$ws.Range("I77").Value = 1800.4286
$ws.Range("K77").Value = 9002.143
$ws.Range("M77").Value = -4634.143

# Row 97
$ws.Range("H97").Value = 727
$ws.Range("I97").Value = 626
$ws.Range("J97").Value = 929
$ws.Range("K97").Value = 626
$ws.Range("L97").Value = 929
$ws.Range("M97").Value = -130
$ws.Range("N97").Value = -1921

# Row 102
$ws.Range("H102").Value = 7885.5713
$ws.Range("I102").Value = 2600
$ws.Range("J102").Value = 9999.799999999999
$ws.Range("K102").Value = 2600
$ws.Range("L102").Value = 9999.799999999999
$ws.Range("M102").Value = -978
$ws.Range("N102").Value = -13243.8

# Row 132
$ws.Range("H132").Value = 2741.861
$ws.Range("I132").Value = 2096.7144
$ws.Range("J132").Value = 4999.875
$ws.Range("K132").Value = 6290.1432
$ws.Range("L132").Value = 14999.625
$ws.Range("M132").Value = -3760.1432
$ws.Range("N132").Value = -20059.625

# Row 136
$ws.Range("H136").Value = 3624.2273
$ws.Range("I136").Value = 2631.4119
$ws.Range("K136").Value = 7894.2357
$ws.Range("M136").Value = -5344.2357

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 12469.5
$ws.Range("I20").Value = 14997
$ws.Range("K20").Value = 14997
$ws.Range("M20").Value = -14750

# Row 134
$ws.Range("H134").Value = 5059.773
$ws.Range("I134").Value = 4824.524
$ws.Range("K134").Value = 14473.572
$ws.Range("M134").Value = -11938.572

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -41250

# Row 60
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022
$ws.Range("M60").ClearContents()

# Row 80
$ws.Range("H80").Value = 86555.39999999999
$ws.Range("J80").Value = 86555.39999999999
$ws.Range("L80").Value = 86555.39999999999
$ws.Range("N80").Value = -88801.39999999999

# Row 83
$ws.Range("H83").Value = 86555.39999999999
$ws.Range("J83").Value = 86555.39999999999
$ws.Range("L83").Value = 259666.2
$ws.Range("N83").Value = -270898.2

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 35.409092
$ws.Range("I2").Value = 31.928572
$ws.Range("J2").Value = 41.5
$ws.Range("K2").Value = 191.571432
$ws.Range("L2").Value = 249
$ws.Range("M2").Value = -78.57143199999999
$ws.Range("N2").Value = -475

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 35740.5
$ws.Range("I132").Value = 42119.215
$ws.Range("K132").Value = 126357.645
$ws.Range("M132").Value = -123827.645

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3144
$ws.Range("I22").Value = 930
$ws.Range("J22").Value = 4251
$ws.Range("K22").Value = 930
$ws.Range("L22").Value = 4251
$ws.Range("M22").Value = -635
$ws.Range("N22").Value = -4841

# Row 27
$ws.Range("H27").Value = 3144
$ws.Range("I27").Value = 930
$ws.Range("J27").Value = 4251
$ws.Range("K27").Value = 930
$ws.Range("L27").Value = 4251
$ws.Range("M27").Value = -823
$ws.Range("N27").Value = -4465

# Row 40
$ws.Range("H40").Value = 6444.3076
$ws.Range("I40").Value = 5798.8184
$ws.Range("K40").Value = 5798.8184
$ws.Range("M40").Value = -5662.8184

# Row 55
$ws.Range("H55").Value = 1302.4
$ws.Range("I55").Value = 1773.8
$ws.Range("J55").Value = 1145.2667
$ws.Range("K55").Value = 1773.8
$ws.Range("L55").Value = 1145.2667
$ws.Range("M55").Value = -1600.8
$ws.Range("N55").Value = -1491.2667

# Row 61
$ws.Range("H61").Value = 3110.7778
$ws.Range("I61").Value = 2260.4783
$ws.Range("K61").Value = 2260.4783
$ws.Range("M61").Value = -2058.4783

# Row 68
$ws.Range("H68").Value = 6542.857
$ws.Range("J68").Value = 8359.799999999999
$ws.Range("L68").Value = 8359.799999999999
$ws.Range("N68").Value = -9857.799999999999

# Row 71
$ws.Range("H71").Value = 6542.857
$ws.Range("J71").Value = 8359.799999999999
$ws.Range("L71").Value = 41799
$ws.Range("N71").Value = -49287

# Row 93
$ws.Range("H93").Value = 2482.4
$ws.Range("I93").Value = 2421
$ws.Range("K93").Value = 2421
$ws.Range("M93").Value = -1173

# Row 96
$ws.Range("H96").Value = 60000
$ws.Range("J96").Value = 60000
$ws.Range("L96").Value = 60000
$ws.Range("N96").Value = -65492

# Row 109
$ws.Range("H109").Value = 59999
$ws.Range("J109").Value = 59999
$ws.Range("L109").Value = 59999
$ws.Range("N109").Value = -62773

# Row 113
$ws.Range("H113").Value = 3110.7778
$ws.Range("I113").Value = 2260.4783
$ws.Range("K113").Value = 2260.4783
$ws.Range("M113").Value = -90.47830000000022

# Row 132
$ws.Range("H132").Value = 2645.3462
$ws.Range("I132").Value = 2249.2354
$ws.Range("J132").Value = 3393.5557
$ws.Range("K132").Value = 6747.706200000001
$ws.Range("L132").Value = 10180.6671
$ws.Range("M132").Value = -4217.706200000001
$ws.Range("N132").Value = -15240.6671

# Row 136
$ws.Range("H136").Value = 4999.5
$ws.Range("I136").Value = 4999.3335
$ws.Range("K136").Value = 14998.0005
$ws.Range("M136").Value = -12448.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10285.571
$ws.Range("I62").Value = 5999.5
$ws.Range("K62").Value = 5999.5
$ws.Range("M62").Value = -5375.5

# Row 65
$ws.Range("H65").Value = 10285.571
$ws.Range("I65").Value = 5999.5
$ws.Range("K65").Value = 29997.5
$ws.Range("M65").Value = -26877.5

# Row 81
$ws.Range("H81").Value = 1099.5
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937

# Row 84
$ws.Range("H84").Value = 1099.5
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686

# Row 100
$ws.Range("H100").Value = 452.86667
$ws.Range("I100").Value = 265.8889
$ws.Range("K100").Value = 531.7778
$ws.Range("M100").Value = 9.222200000000043

# Row 126
$ws.Range("H126").Value = 5419.4546
$ws.Range("I126").Value = 3269
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 9807
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -7337
$ws.Range("N126").Value = -28940

# Row 132
$ws.Range("H132").Value = 3594.7307
$ws.Range("I132").Value = 3366.5789
$ws.Range("K132").Value = 10099.7367
$ws.Range("M132").Value = -7569.736699999999

# Row 136
$ws.Range("H136").Value = 4995.609
$ws.Range("I136").Value = 3247.75
$ws.Range("K136").Value = 9743.25
$ws.Range("M136").Value = -7193.25
